$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "fatma"
$ws.Range("B2").Value = "2024-11-11 11:00"

$ws.Range("A3").Value = "fatma"
$ws.Range("B3").Value = "2024-11-11 12:00"

$ws.Range("A4").Value = "ayşe"
$ws.Range("B4").Value = "2024-11-11 13:00"

$ws.Range("A5").Value = "ayşe"
$ws.Range("B5").Value = "2024-11-11 14:00"
